$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of OHLCV candle data to append (rows 1169-1232)
$csvData = @"
1169,45534.5,1.435,1.455,1.37,1.384,614413.05
1170,45534.58333333334,1.384,1.393,1.325,1.335,698258.8
1171,45534.66666666666,1.335,1.36,1.324,1.352,388242.11
1172,45534.75,1.351,1.408,1.35,1.366,291343.17
1173,45534.83333333334,1.366,1.426,1.364,1.398,332550.67
1174,45534.91666666666,1.4,1.427,1.395,1.426,157487.12
1175,45535,1.426,1.509,1.414,1.494,913033.91
1176,45535.08333333334,1.494,1.5,1.464,1.467,263271.94
1177,45535.16666666666,1.467,1.48,1.465,1.468,88017.06
1178,45535.25,1.468,1.49,1.457,1.475,363039.38
1179,45535.33333333334,1.476,1.487,1.451,1.455,268251.64
1180,45535.41666666666,1.454,1.477,1.45,1.465,202277.67
1181,45535.5,1.465,1.465,1.4,1.404,905602.39
1182,45535.58333333334,1.403,1.403,1.366,1.371,378311.6
1183,45535.66666666666,1.372,1.372,1.337,1.344,227071.77
1184,45535.75,1.343,1.358,1.338,1.352,217073.7
1185,45535.83333333334,1.352,1.368,1.343,1.361,140031.81
1186,45535.91666666666,1.365,1.379,1.36,1.363,295107.05
1187,45536,1.362,1.363,1.343,1.344,219387.96
1188,45536.08333333334,1.344,1.348,1.336,1.338,148364.92
1189,45536.16666666666,1.338,1.349,1.325,1.327,207835.56
1190,45536.25,1.328,1.338,1.324,1.335,229726.92
1191,45536.33333333334,1.334,1.335,1.314,1.317,264602.71
1192,45536.41666666666,1.317,1.323,1.304,1.318,269940.49
1193,45536.5,1.317,1.319,1.278,1.298,293251.7
1194,45536.58333333334,1.299,1.324,1.276,1.319,324983.64
1195,45536.66666666666,1.319,1.371,1.311,1.315,802923.4
1196,45536.75,1.315,1.355,1.309,1.334,311398.56
1197,45536.83333333334,1.334,1.346,1.3,1.317,238594.99
1198,45536.91666666666,1.316,1.316,1.257,1.287,867566.52
1199,45537,1.286,1.297,1.263,1.272,283269.16
1200,45537.08333333334,1.272,1.292,1.261,1.281,154975.21
1201,45537.16666666666,1.281,1.282,1.248,1.262,653690.59
1202,45537.25,1.263,1.273,1.243,1.246,274714.42
1203,45537.33333333334,1.245,1.288,1.238,1.285,375849.98
1204,45537.41666666666,1.284,1.293,1.26,1.27,629702.54
1205,45537.5,1.269,1.278,1.246,1.25,193008.65
1206,45537.58333333334,1.25,1.271,1.236,1.249,258263.64
1207,45537.66666666666,1.249,1.266,1.236,1.264,333901.63
1208,45537.75,1.263,1.266,1.247,1.248,145124.95
1209,45537.83333333334,1.248,1.257,1.246,1.255,18111.08
1210,45537.91666666666,1.268,1.28,1.266,1.268,290666.84
1211,45538,1.268,1.275,1.256,1.27,174719.46
1212,45538.08333333334,1.27,1.287,1.254,1.257,323699.14
1213,45538.16666666666,1.258,1.264,1.249,1.263,107821.78
1214,45538.25,1.263,1.263,1.244,1.253,111636.69
1215,45538.33333333334,1.253,1.254,1.223,1.225,344852.99
1216,45538.41666666666,1.225,1.236,1.222,1.225,165382.12
1217,45538.5,1.225,1.236,1.193,1.199,529652.6
1218,45538.58333333334,1.198,1.199,1.171,1.174,535505.21
1219,45538.66666666666,1.174,1.19,1.17,1.18,234686.59
1220,45538.75,1.181,1.2,1.178,1.188,409774.21
1221,45538.83333333334,1.187,1.193,1.179,1.183,177022.72
1222,45538.91666666666,1.183,1.186,1.169,1.172,187001.17
1223,45539,1.172,1.187,1.121,1.162,895197.28
1224,45539.08333333334,1.161,1.187,1.157,1.183,200246.25
1225,45539.16666666666,1.183,1.19,1.177,1.177,209103.31
1226,45539.25,1.177,1.209,1.177,1.207,576232.26
1227,45539.33333333334,1.208,1.208,1.2,1.208,35556.85
1228,45539.41666666666,1.216,1.232,1.197,1.199,229153.54
1229,45539.5,1.199,1.219,1.196,1.215,239468.67
1230,45539.58333333334,1.215,1.24,1.203,1.239,399959.58
1231,45539.66666666666,1.239,1.311,1.236,1.243,1621284.68
1232,45539.75,1.243,1.249,1.233,1.238,136020.2
"@

# Apply the existing date/number style (used by column A, e.g. A1168) to the
# new column-A cells so they keep the "YYYY-MM-DD HH:MM:SS" style (s="2").
$ws.Range("A1168").Copy()
$ws.Range("A1169:A1232").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($line in ($csvData -split "`r?`n")) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    $ws.Cells.Item($row, 1).Value = [double]$parts[1]
    $ws.Cells.Item($row, 2).Value = [double]$parts[2]
    $ws.Cells.Item($row, 3).Value = [double]$parts[3]
    $ws.Cells.Item($row, 4).Value = [double]$parts[4]
    $ws.Cells.Item($row, 5).Value = [double]$parts[5]
    $ws.Cells.Item($row, 6).Value = [double]$parts[6]
}
